# "Lista de materiales CNC Laser" — add the "ENLACES - OPCION 2" column and
# the two extra laser rows described in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column before the old "NOTAS" column (old E -> F) ---
$ws.Columns("E:E").Insert()

# Match the new column's width to column D's width (they end up the same
# width in the source workbook).
$ws.Columns("E:E").ColumnWidth = $ws.Columns("D:D").ColumnWidth

# --- 2. New column header + "alternate link" values ---
$ws.Range("E3").Value2  = "ENLACES - OPCION 2"

$ws.Range("E5").Value2  = "https://www.banggood.com/custlink/GKvmsnZmI6"
$ws.Range("E10").Value2 = "https://www.banggood.com/custlink/3mmKsnZKSf"
$ws.Range("E11").Value2 = "https://www.banggood.com/custlink/KGvKAQivaL"
$ws.Range("E14").Value2 = "https://www.banggood.com/custlink/3D33NbiDO1"
$ws.Range("E15").Value2 = "https://www.banggood.com/custlink/mvDGsbi3gR"
$ws.Range("E16").Value2 = "https://www.banggood.com/custlink/GmKmssZv6s"
$ws.Range("E26").Value2 = "https://www.banggood.com/custlink/KvDvAb9DsG"
$ws.Range("E28").Value2 = "https://www.banggood.com/custlink/mGDGnAZ3Ac"
$ws.Range("E29").Value2 = "https://www.banggood.com/custlink/mDvGsnIDAg"

# --- 3. New rows with info/links about the laser module options ---
$ws.Range("E35").Value2 = "LASER 7W PWM"
$ws.Hyperlinks.Add($ws.Range("F35"), "https://www.banggood.com/custlink/mDGvQAZDkK") | Out-Null

$ws.Range("F36").Value2 = "https://www.banggood.com/custlink/vm3GNnIKkJ"
$ws.Range("D36").Value2 = " "
$ws.Range("E36").Value2 = "LASER 3,5W PWM - es el que tengo en casa"

$ws.Rows("35:35").RowHeight = 15.75
$ws.Rows("36:36").RowHeight = 15.75

# --- 4. Move the active selection, matching the author's saved cursor ---
$ws.Range("E39").Select()
